# Updated cryptos list with latest price/volume data (per GitHub Actions scrape)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D (Price) and E (Volume) cells hold plain text in the source data;
# force text format so numeric-looking strings (e.g. '27.10', '1.00') keep
# their exact formatting instead of being coerced to numbers.
$cells = @("B45","C45","B46","C46")
$cells += "D2"
$cells += "E2"
$cells += "D3"
$cells += "E3"
$cells += "E4"
$cells += "D5"
$cells += "E5"
$cells += "D6"
$cells += "E6"
$cells += "E7"
$cells += "D8"
$cells += "D9"
$cells += "E9"
$cells += "D10"
$cells += "E10"
$cells += "E11"
$cells += "E12"
$cells += "E13"
$cells += "D14"
$cells += "E14"
$cells += "E15"
$cells += "D16"
$cells += "E16"
$cells += "D17"
$cells += "D18"
$cells += "E18"
$cells += "E19"
$cells += "D20"
$cells += "E20"
$cells += "D21"
$cells += "E21"
$cells += "E22"
$cells += "D23"
$cells += "E23"
$cells += "E24"
$cells += "D25"
$cells += "E25"
$cells += "D26"
$cells += "E26"
$cells += "D27"
$cells += "E27"
$cells += "D28"
$cells += "E28"
$cells += "D29"
$cells += "E29"
$cells += "E30"
$cells += "E31"
$cells += "D32"
$cells += "E32"
$cells += "E33"
$cells += "D34"
$cells += "E34"
$cells += "D35"
$cells += "E35"
$cells += "E36"
$cells += "D37"
$cells += "E37"
$cells += "E38"
$cells += "D39"
$cells += "E39"
$cells += "D40"
$cells += "E40"
$cells += "D41"
$cells += "E41"
$cells += "D42"
$cells += "E42"
$cells += "D43"
$cells += "E43"
$cells += "D45"
$cells += "E45"
$cells += "D46"
$cells += "E46"
$cells += "D47"
$cells += "E47"
$cells += "D48"
$cells += "E48"
$cells += "D49"
$cells += "E49"
$cells += "E50"
$cells += "D51"
$cells += "E51"
foreach ($addr in $cells) {
    $ws.Range($addr).NumberFormat = "@"
}

$ws.Range("D2").Value = '62.656.84'
$ws.Range("E2").Value = '  -1.54%  '

$ws.Range("D3").Value = '2.443.93'
$ws.Range("E3").Value = '  -1.77%  '

$ws.Range("E4").Value = '  -0.32%  '

$ws.Range("D5").Value = '569.17'
$ws.Range("E5").Value = '  -0.86%  '

$ws.Range("D6").Value = '144.91'
$ws.Range("E6").Value = '  -3.19%  '

$ws.Range("E7").Value = '  +0.20%  '

$ws.Range("D8").Value = '0.532'

$ws.Range("D9").Value = '2.439.89'
$ws.Range("E9").Value = '  -2.42%  '

$ws.Range("D10").Value = '0.109'
$ws.Range("E10").Value = '  -4.11%  '

$ws.Range("E11").Value = '  +1.32%  '

$ws.Range("E12").Value = '  -2.43%  '

$ws.Range("E13").Value = '  -2.71%  '

$ws.Range("D14").Value = '27.10'
$ws.Range("E14").Value = '  -0.69%  '

$ws.Range("E15").Value = '  -5.03%  '

$ws.Range("D16").Value = '2.887.83'
$ws.Range("E16").Value = '  -1.47%  '

$ws.Range("D17").Value = '62.461.05'

$ws.Range("D18").Value = '2.447.62'
$ws.Range("E18").Value = '  -2.08%  '

$ws.Range("E19").Value = '  -2.83%  '

$ws.Range("D20").Value = '7.22'
$ws.Range("E20").Value = '  -1.04%  '

$ws.Range("D21").Value = '327.38'
$ws.Range("E21").Value = '  -0.56%  '

$ws.Range("E22").Value = '  -1.89%  '

$ws.Range("D23").Value = '2.09'
$ws.Range("E23").Value = '  +10.38%  '

$ws.Range("E24").Value = '  -0.35%  '

$ws.Range("D25").Value = '65.16'
$ws.Range("E25").Value = '  -3.55%  '

$ws.Range("D26").Value = '625.41'
$ws.Range("E26").Value = '  -3.00%  '

$ws.Range("D27").Value = '9.02'
$ws.Range("E27").Value = '  +2.51%  '

$ws.Range("D28").Value = '0.0₃0995'
$ws.Range("E28").Value = '  -6.03%  '

$ws.Range("D29").Value = '2.556.96'
$ws.Range("E29").Value = '  -2.15%  '

$ws.Range("E30").Value = '  -2.12%  '

$ws.Range("E31").Value = '  +0.66%  '

$ws.Range("D32").Value = '8.14'
$ws.Range("E32").Value = '  -4.87%  '

$ws.Range("E33").Value = '  -2.61%  '

$ws.Range("D34").Value = '0.138'
$ws.Range("E34").Value = '  -4.64%  '

$ws.Range("D35").Value = '5.12'
$ws.Range("E35").Value = '  -1.21%  '

$ws.Range("E36").Value = '  -3.50%  '

$ws.Range("D37").Value = '1.00'
$ws.Range("E37").Value = '  +0.29%  '

$ws.Range("E38").Value = '  -2.91%  '

$ws.Range("D39").Value = '18.82'
$ws.Range("E39").Value = '  -0.91%  '

$ws.Range("D40").Value = '5.30'
$ws.Range("E40").Value = '  -4.69%  '

$ws.Range("D41").Value = '146.59'
$ws.Range("E41").Value = '  -0.69%  '

$ws.Range("D42").Value = '1.77'
$ws.Range("E42").Value = '  -4.70%  '

$ws.Range("D43").Value = '2.57'
$ws.Range("E43").Value = '  -3.27%  '

$ws.Range("B45").Value = 'Aave'
$ws.Range("C45").Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range("D45").Value = '146.22'
$ws.Range("E45").Value = '  -4.21%  '

$ws.Range("B46").Value = 'Filecoin'
$ws.Range("C46").Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range("D46").Value = '3.75'
$ws.Range("E46").Value = '  -0.96%  '

$ws.Range("D47").Value = '20.67'
$ws.Range("E47").Value = '  -2.41%  '

$ws.Range("D48").Value = '0.0530'
$ws.Range("E48").Value = '  -3.97%  '

$ws.Range("D49").Value = '0.596'
$ws.Range("E49").Value = '  -3.01%  '

$ws.Range("E50").Value = '  -3.12%  '

$ws.Range("D51").Value = '0.0921'
$ws.Range("E51").Value = '  -1.10%  '
